$wb = $excel.ActiveWorkbook

# OFF sheet - Week 17 "Road" row (row 3) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 209
$wsOff.Range("C3").Value = 153
$wsOff.Range("D3").Value = 53
$wsOff.Range("E3").Value = 23

# DEF sheet - Week 17 "Road" row (row 3) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 206
$wsDef.Range("C3").Value = 130
$wsDef.Range("D3").Value = 50
$wsDef.Range("E3").Value = 26
